$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4:D4").Value = "Manage Memberships"
$ws.Range("A5:D5").Value = "Manage Members"
